$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.300.05"
$ws.Range("E2").Value = "  +3.84%  "
$ws.Range("E3").Value = "  +4.67%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'328.91"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4462"
$ws.Range("E7").Value = "  +5.38%  "
$ws.Range("E8").Value = "  +3.43%  "
$ws.Range("D9").Value = "'44.91"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.07707"
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "'22.03"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").Value = "'6.300"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "'7.573"
$ws.Range("E15").Value = "  +6.09%  "
$ws.Range("D16").Value = "1.845.19"
$ws.Range("E16").Value = "  +6.63%  "
$ws.Range("D17").Value = "'93.10"
$ws.Range("E17").Value = "  +7.15%  "
$ws.Range("D18").Value = "'0.00001082"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").Value = "'0.06534"
$ws.Range("E19").Value = "  +9.43%  "
$ws.Range("D20").Value = "'0.9998"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'17.52"
$ws.Range("E21").Value = "  +4.06%  "
$ws.Range("D22").Value = "'6.227"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").Value = "28.326.76"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("D24").Value = "'11.68"
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("D25").Value = "'2.170"
$ws.Range("E25").Value = "  -8.88%  "
$ws.Range("D26").Value = "'20.79"
$ws.Range("E26").Value = "  +3.57%  "
$ws.Range("D27").Value = "'155.83"
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("D28").Value = "2.042.15"
$ws.Range("E28").Value = "  +6.03%  "
$ws.Range("D29").Value = "'2.310"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").Value = "'128.43"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Value = "'1.199"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "'5.910"
$ws.Range("E32").Value = "  +5.76%  "
$ws.Range("D33").Value = "'0.09253"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").Value = "'13.03"
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("D36").Value = "'0.02354"
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("D37").Value = "'0.2184"
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("D38").Value = "'5.179"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("D39").Value = "'0.06230"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("D40").Value = "'0.6574"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").Value = "'1.198"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "'8.145"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "'1.408"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "'13.95"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").Value = "'0.6087"
$ws.Range("E46").Value = "  +4.84%  "
$ws.Range("D47").Value = "'3.766"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "'126.99"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").Value = "'2.033"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("D50").Value = "'1.157"
$ws.Range("E50").Value = "  +5.81%  "
$ws.Range("E51").Value = "  +2.60%  "
